$p = $ppt.ActivePresentation

# --- Step 1: duplicate the last slide ("Recap") and move the duplicate ---
# before it, so the deck ends up with:
#   position 9  -> new slide (SlideID 265) "Recap" (unchanged content)
#   position 10 -> original slide (SlideID 264) repurposed as
#                  "Pros and Cons of Django"
$recap = $p.Slides.Item($p.Slides.Count)
$dup = $recap.Duplicate()
$dup.MoveTo($recap.SlideIndex)

# --- Step 2: rewrite the original slide (now last) with the new content ---
$target = $p.Slides.Item($p.Slides.Count)

$title = $target.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Pros and Cons of Django"

$body = $target.Shapes.Item(2)
$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = "Pros`rDjango is a " + [char]0x201C + "High Level" + [char]0x201D + " framework i.e., it offers a lot of inbuilt tools and utilities e.g., admin panel, user authentication or testing-libraries`rScalability - Django is built to handle millions of users`rClearly structured (MVC)`rLess code because of reusable apps`rCons`rDjango is known for taking up a lot of resources (less suited for small projects)`rSlower compared to other backend frameworks`r"

$bodyTr.Paragraphs(1).IndentLevel = 2
$bodyTr.Paragraphs(2).IndentLevel = 3
$bodyTr.Paragraphs(3).IndentLevel = 3
$bodyTr.Paragraphs(4).IndentLevel = 3
$bodyTr.Paragraphs(5).IndentLevel = 3
$bodyTr.Paragraphs(6).IndentLevel = 2
$bodyTr.Paragraphs(7).IndentLevel = 3
$bodyTr.Paragraphs(8).IndentLevel = 3
$bodyTr.Paragraphs(9).IndentLevel = 3

# --- Step 3: refresh the cached "today" date field shown in the footers ---
# (slide master + every slide layout), from 06.11.2021 to 11.11.2021.
$newDate = "11.11.2021"
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "06.11.2021") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "06.11.2021") {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
